$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.07205133333333334
$ws.Range("M2").Value = 1.484826
$ws.Range("N2").Value = 4.454478
$ws.Range("O2").Value = 0.06049021884829667
$ws.Range("P2").Value = 0.06049021884829667
$ws.Range("Q2").Value = 0.106983693068
$ws.Range("R2").Value = 0.962853237612
$ws.Range("S2").Value = 0.06049021884829667
$ws.Range("T2").Value = 0.06049021884829667

# Row 3
$ws.Range("G3").Value = 0.07205133333333334
$ws.Range("O3").Value = 0.5859425360316464
$ws.Range("P3").Value = 0.5859425360316464
$ws.Range("R3").Value = 9.326742052092001
$ws.Range("S3").Value = 0.5859425360316464
$ws.Range("T3").Value = 0.5859425360316464

# Row 4
$ws.Range("G4").Value = 0.07205133333333334
$ws.Range("M4").Value = 8.653369666666666
$ws.Range("N4").Value = 25.960109
$ws.Range("O4").Value = 0.3525289999716321
$ws.Range("P4").Value = 0.3525289999716321
$ws.Range("Q4").Value = 0.6234868223095557
$ws.Range("R4").Value = 5.611381400786001
$ws.Range("S4").Value = 0.3525289999716321
$ws.Range("T4").Value = 0.3525289999716321

# Row 5
$ws.Range("G5").Value = 0.07205133333333334
$ws.Range("M5").Value = 0.02548533333333333
$ws.Range("N5").Value = 0.076456
$ws.Range("O5").Value = 0.001038245148424882
$ws.Range("P5").Value = 0.001038245148424882
$ws.Range("Q5").Value = 0.001836252247111111
$ws.Range("R5").Value = 0.016526270224
$ws.Range("S5").Value = 0.001038245148424882
$ws.Range("T5").Value = 0.001038245148424882
